# Doctor record export (#58)
# The "Anesthesiology_Physician" column header on the first sheet was
# renamed to "Anesthesiologist_Physician" (matching the spelling already
# used on the second sheet), and the workbook's active/selected tab moved
# from the second sheet back to the first sheet.

$wb = $excel.ActiveWorkbook

# First worksheet: "Nov 19 2020 - Dec 1 2020"
$ws1 = $wb.Worksheets.Item(1)

# Fix the column I header spelling.
$ws1.Range("I1").Value = "Anesthesiologist_Physician"

# Make the first sheet the active / selected sheet (tab), matching the
# target workbook view state.
$ws1.Activate()
